# Update column F (dSF) values on Sheet1 to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F11").Value = -2
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = -10
$ws.Range("F35").Value = -2
$ws.Range("F42").Value = -2
$ws.Range("F46").Value = 0
